$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1031.1666
$ws.Range("I2").Value = 1031.1666
$ws.Range("K2").Value = 1031.1666
$ws.Range("M2").Value = -918.1666
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H49").Value = 249
$ws.Range("I49").Value = 199
$ws.Range("K49").Value = 597
$ws.Range("M49").Value = -461
$ws.Range("H107").Value = 494.33334
$ws.Range("J107").Value = 99.5
$ws.Range("L107").Value = 99.5
$ws.Range("N107").Value = -3939.5
$ws.Range("H132").Value = 825
$ws.Range("I132").Value = 820.7317
$ws.Range("K132").Value = 2462.1951
$ws.Range("M132").Value = 67.80490000000009
$ws.Range("H137").Value = 1823.9395
$ws.Range("I137").Value = 1457.1818
$ws.Range("K137").Value = 4371.5454
$ws.Range("M137").Value = -1821.5454
$ws.Range("H138").Value = 4306.069
$ws.Range("I138").Value = 3749
$ws.Range("J138").Value = 4395.2
$ws.Range("K138").Value = 11247
$ws.Range("L138").Value = 13185.6
$ws.Range("M138").Value = -6107
$ws.Range("N138").Value = -23465.6
$ws.Range("H141").Value = 8666.166999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 283
$ws.Range("I97").Value = 283
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 283
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 213
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 33500
$ws.Range("J95").Value = 33500
$ws.Range("L95").Value = 33500
$ws.Range("N95").Value = -38992
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("N97").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3045.4
$ws.Range("I31").Value = 1876
$ws.Range("J31").Value = 3825
$ws.Range("K31").Value = 1876
$ws.Range("L31").Value = 3825
$ws.Range("M31").Value = -1581
$ws.Range("N31").Value = -4415
$ws.Range("H34").Value = 3045.4
$ws.Range("I34").Value = 1876
$ws.Range("J34").Value = 3825
$ws.Range("K34").Value = 1876
$ws.Range("L34").Value = 3825
$ws.Range("M34").Value = -1674
$ws.Range("N34").Value = -4229
$ws.Range("H86").Value = 22538.555
$ws.Range("J86").Value = 31907
$ws.Range("L86").Value = 31907
$ws.Range("N86").Value = -34153
$ws.Range("H89").Value = 22538.555
$ws.Range("J89").Value = 31907
$ws.Range("L89").Value = 159535
$ws.Range("N89").Value = -170767

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 350
$ws.Range("I10").Value = 350
$ws.Range("K10").Value = 1050
$ws.Range("M10").Value = -911
$ws.Range("H25").Value = 599.5
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 599.5
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = 1798.5
$ws.Range("N25").Value = -2136.5
$ws.Range("H30").Value = 599.5
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 599.5
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("M30").Value = 1798.5
$ws.Range("N30").Value = -2002.5
$ws.Range("H34").Value = 3519.6
$ws.Range("I34").Value = 1366
$ws.Range("J34").Value = 6750
$ws.Range("K34").Value = 4098
$ws.Range("L34").Value = 20250
$ws.Range("M34").Value = -4014
$ws.Range("N34").Value = -20418
$ws.Range("H46").Value = 2850
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H54").Value = 2066
$ws.Range("J54").Value = 2066
$ws.Range("L54").Value = 6198
$ws.Range("N54").Value = -7316
$ws.Range("H55").Value = 6738.8
$ws.Range("I55").Value = 1398
$ws.Range("J55").Value = 14750
$ws.Range("K55").Value = 4194
$ws.Range("L55").Value = 44250
$ws.Range("M55").Value = -4017
$ws.Range("N55").Value = -44604
$ws.Range("H56").Value = 19608.666
$ws.Range("I56").Value = 19608.666
$ws.Range("K56").Value = 19608.666
$ws.Range("M56").Value = -19078.666
$ws.Range("H113").Value = 4067.5715
$ws.Range("J113").Value = 4067.5715
$ws.Range("L113").Value = 12202.7145
$ws.Range("N113").Value = -16542.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H82").Value = 2154.2
$ws.Range("I82").Value = 2269.6667
$ws.Range("K82").Value = 2269.6667
$ws.Range("M82").Value = -1908.6667
$ws.Range("H85").Value = 2154.2
$ws.Range("I85").Value = 2269.6667
$ws.Range("K85").Value = 2269.6667
$ws.Range("M85").Value = -1021.6667
$ws.Range("H100").Value = 1749.5
$ws.Range("I100").Value = 499
$ws.Range("K100").Value = 499
$ws.Range("M100").Value = 42
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 3559.818
$ws.Range("I136").Value = 3727.8518
$ws.Range("K136").Value = 11183.5554
$ws.Range("M136").Value = -8633.555399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4394.706
$ws.Range("I81").Value = 2615.4285
$ws.Range("J81").Value = 12698
$ws.Range("K81").Value = 5230.857
$ws.Range("L81").Value = 25396
$ws.Range("M81").Value = -4169.857
$ws.Range("N81").Value = -27518
$ws.Range("H84").Value = 4394.706
$ws.Range("I84").Value = 2615.4285
$ws.Range("J84").Value = 12698
$ws.Range("K84").Value = 26154.285
$ws.Range("L84").Value = 126980
$ws.Range("M84").Value = -20850.285
$ws.Range("N84").Value = -137588
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H100").Value = 1172.5714
$ws.Range("I100").Value = 1172.5714
$ws.Range("K100").Value = 2345.1428
$ws.Range("M100").Value = -1804.1428
$ws.Range("H122").Value = 4474.5
$ws.Range("I122").Value = 4474.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13423.5
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10973.5
$ws.Range("H126").Value = 2111
$ws.Range("I126").Value = 2051.7334
$ws.Range("K126").Value = 6155.2002
$ws.Range("M126").Value = -3685.2002
